$d = $word.ActiveDocument

# Paragraph 5 ("<id>...</id>") currently holds three separate runs:
#   1. "<id>"     - Courier New, 7f6000, 18pt  (the XML-tag styling)
#   2. "p056v_1"  - plain black                (the id value)
#   3. "</id>"    - Courier New, 7f6000, 18pt  (the XML-tag styling)
# The edit collapses these into a single run reading "<id>p056v_1</id>"
# that keeps the tag run's formatting (runs 1 and 3's rPr).

# Locate the middle ("value") run and remember where it starts.
$valueRng = $d.Content
$found = $valueRng.Find.Execute("p056v_1", $true, $false, $false, $false, `
                                 $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "could not find the '<id>' value run (p056v_1)"
}
$insertAt = $valueRng.Start

# Deleting it leaves the two neighboring "<id>" / "</id>" runs directly
# adjacent; since they carry identical character formatting Word merges
# them into one run reading "<id></id>".
$valueRng.Delete()

# Insert the id value back at the same offset. The insertion point now
# sits inside the merged tag run, so the inserted text inherits that
# run's formatting, yielding a single run: "<id>p056v_1</id>".
$insertionPoint = $d.Range($insertAt, $insertAt)
$insertionPoint.InsertAfter("p056v_1")
